$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D (Price) to remain text so numeric-looking values
# ("42.207.27", "1.01", etc.) are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.207.27"
$ws.Range("E2").Value = "  -1.17%  "
$ws.Range("D3").Value = "2.304.60"
$ws.Range("E3").Value = "  -2.28%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "312.43"
$ws.Range("E5").Value = "  -3.98%  "
$ws.Range("D6").Value = "105.48"
$ws.Range("E6").Value = "  +4.39%  "
$ws.Range("D7").Value = "0.627"
$ws.Range("E7").Value = "  -1.59%  "
$ws.Range("E8").Value = "  +0.16%  "
$ws.Range("D9").Value = "0.608"
$ws.Range("E9").Value = "  -2.69%  "
$ws.Range("D10").Value = "40.36"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").Value = "0.0914"
$ws.Range("E11").Value = "  -0.88%  "
$ws.Range("D12").Value = "8.28"
$ws.Range("E12").Value = "  -2.13%  "
$ws.Range("E13").Value = "  -0.06%  "
$ws.Range("D14").Value = "0.975"
$ws.Range("E14").Value = "  -2.84%  "
$ws.Range("D15").Value = "15.57"
$ws.Range("E15").Value = "  -5.74%  "
$ws.Range("D16").Value = "2.655.22"
$ws.Range("E16").Value = "  -2.12%  "
$ws.Range("D17").Value = "2.306.94"
$ws.Range("E17").Value = "  -1.93%  "
$ws.Range("D18").Value = "42.131.31"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "7.64"
$ws.Range("E19").Value = "  -5.15%  "
$ws.Range("E20").Value = "  -1.45%  "
$ws.Range("D21").Value = "74.63"
$ws.Range("E21").Value = "  -2.05%  "
$ws.Range("D22").Value = "3.48"
$ws.Range("E22").Value = "  -6.13%  "
$ws.Range("D23").Value = "259.02"
$ws.Range("E23").Value = "  -2.71%  "
$ws.Range("D25").Value = "9.33"
$ws.Range("E25").Value = "  -7.41%  "
$ws.Range("D26").Value = "1.01"
$ws.Range("E26").Value = "  +0.38%  "
$ws.Range("D27").Value = "10.99"
$ws.Range("E27").Value = "  -4.18%  "
$ws.Range("D29").Value = "22.81"
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("D30").Value = "35.83"
$ws.Range("E30").Value = "  +1.49%  "
$ws.Range("D31").Value = "164.85"
$ws.Range("E31").Value = "  -6.29%  "
$ws.Range("D32").Value = "0.0899"
$ws.Range("E32").Value = "  +0.00%  "
$ws.Range("E33").Value = "  -5.51%  "
$ws.Range("D34").Value = "5.85"
$ws.Range("E34").Value = "  -3.33%  "
$ws.Range("B35").Value = "Kaspa"
$ws.Range("C35").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D35").Value = "0.119"
$ws.Range("E35").Value = "  +12.25%  "
$ws.Range("B36").Value = "Stellar"
$ws.Range("C36").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D36").Value = "0.130"
$ws.Range("E36").Value = "  -1.78%  "
$ws.Range("D37").Value = "4.56"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "0.0354"
$ws.Range("E38").Value = "  -1.23%  "
$ws.Range("E39").Value = "  -5.46%  "
$ws.Range("D40").Value = "3.62"
$ws.Range("E40").Value = "  -4.77%  "
$ws.Range("D41").Value = "72.16"
$ws.Range("E41").Value = "  +3.05%  "
$ws.Range("D42").Value = "98.67"
$ws.Range("E42").Value = "  +8.47%  "
$ws.Range("E43").Value = "  -3.26%  "
$ws.Range("D44").Value = "0.229"
$ws.Range("E44").Value = "  -2.71%  "
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "12.37"
$ws.Range("E46").Value = "  +3.89%  "
$ws.Range("D47").Value = "112.84"
$ws.Range("E47").Value = "  -5.59%  "
$ws.Range("D48").Value = "9.03"
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("D49").Value = "5.33"
$ws.Range("D50").Value = "74.18"
$ws.Range("E50").Value = "  +1.91%  "
$ws.Range("D51").Value = "1.570.37"
$ws.Range("E51").Value = "  +0.77%  "
